# Automatische test-sync: 2025-08-01 23:45:50
# Appends a new log entry (row 10) to the "Logs" sheet and bumps the
# matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 10

$logs.Cells.Item($newRow, 1).Value = "Wil je deze klant bellen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #5: Wil je deze klant bellen?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Beste,
Bedankt voor uw bericht. Om u beter van dienst te kunnen zijn, zou ik graag wat meer informatie willen ontvangen. Kunt u mij alstublieft de naam van de klant en het telefoonnummer doorgeven, zodat wij contact met hen kunnen opnemen?
Met vriendelijke groet,
[Jouw Naam]
E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-08-01 23:45:13"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# The multi-line Antwoord text triggers an automatic row-height change;
# put the row back to its normal (default) auto-fit height so it does
# not pick up an explicit ht/customHeight override.
$logs.Rows.Item($newRow).AutoFit()

# Extend conditional formatting ranges that covered rows 2:9 down to 2:10.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "9")
    $newRange = $logs.Range($col + "2:" + $col + "10")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 2
